$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6897621750831604
$ws.Range("B1").Value = 1.036483645439148
$ws.Range("C1").Value = 0.997683048248291
$ws.Range("D1").Value = 3.945881128311157
$ws.Range("E1").Value = 1.78892970085144
